# Apply updated dSF ("F" column) values to Sheet1.
# These reflect a repull/recalculation of the data (see commit message:
# "repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    3  = -3
    6  = -3
    8  = 0
    13 = 0
    18 = -1
    25 = -3
    30 = -1
    31 = -1
    35 = -3
    36 = -2
    47 = -2
    51 = -2
    54 = -3
    55 = 8
    59 = 0
    64 = -1
    66 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
